$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: height grows (no longer "custom" height target value) ---
$ws.Rows.Item(3).RowHeight = 157.5

# --- New shared strings used by the new rows ---
# (Excel will append these to sharedStrings.xml automatically the first
# time they are used as text cell values.)
$inputFixing = "Input fixing"
$installingNetCode = "Installing Net Code"

# --- New rows 26-44: a simple day-by-day checklist ---
# Column A = day index (1..18), Column B = date (formatted d-mmm),
# Column C = task text for the first two days only.
$startSerial = 45342

for ($i = 0; $i -lt 18; $i++) {
    $row = 26 + $i
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $startSerial + $i
    $ws.Cells.Item($row, 2).NumberFormat = "d-mmm"
    $ws.Rows.Item($row).RowHeight = 50.1
}

$ws.Range("C26").Value = $inputFixing
$ws.Range("C27").Value = $installingNetCode

# Row 44: an extra (empty) row below the list, formatted like the date column.
$ws.Range("A44").NumberFormat = "d-mmm"
$ws.Rows.Item(44).RowHeight = 50.1

# --- Update selection / scroll position to match where the new data is ---
$ws.Range("C27").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
